$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark the "* A an employee, I can view a log of all transactions." task as DONE
$ws.Range("C12").Value = "DONE"

# Select C13 to match the final cursor position
$ws.Range("C13").Select()
